$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1376.2727
$ws.Range("I15").Value = 1376.2727
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4128.8181
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -3959.8181

# Row 99
$ws.Range("H99").Value = 499.5
$ws.Range("I99").Value = 499.33334
$ws.Range("J99").Value = 500
$ws.Range("K99").Value = 1498.00002
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -0.00001999999994950485

# Row 132
$ws.Range("H132").Value = 4142.116
$ws.Range("I132").Value = 4106.8335
$ws.Range("J132").Value = 4323.5713
$ws.Range("K132").Value = 12320.5005
$ws.Range("L132").Value = 12970.7139
$ws.Range("M132").Value = -9790.500499999998

$ws = $wb.Worksheets.Item("ARM")
# Row 80
$ws.Range("H80").Value = 48571
$ws.Range("I80").Value = 47499.75
$ws.Range("J80").Value = 49999.332
$ws.Range("K80").Value = 47499.75
$ws.Range("L80").Value = 49999.332
$ws.Range("M80").Value = -46501.75
$ws.Range("N80").Value = -51995.332

# Row 83
$ws.Range("H83").Value = 48571
$ws.Range("I83").Value = 47499.75
$ws.Range("J83").Value = 49999.332
$ws.Range("K83").Value = 142499.25
$ws.Range("L83").Value = 149997.996
$ws.Range("M83").Value = -137507.25
$ws.Range("N83").Value = -159981.996

$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 5498
$ws.Range("I16").Value = 5498
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5498
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -5328
$ws.Range("N16").ClearContents()

# Row 86
$ws.Range("H86").Value = 874.36365
$ws.Range("I86").Value = 749.9474
$ws.Range("J86").Value = 1662.3334
$ws.Range("K86").Value = 749.9474
$ws.Range("L86").Value = 1662.3334
$ws.Range("M86").Value = 373.0526

# Row 89
$ws.Range("H89").Value = 874.36365
$ws.Range("I89").Value = 749.9474
$ws.Range("J89").Value = 1662.3334
$ws.Range("K89").Value = 3749.737
$ws.Range("L89").Value = 8311.666999999999
$ws.Range("M89").Value = 1866.263

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 770579.9399999999
$ws.Range("I58").Value = 1131.375
$ws.Range("J58").Value = 2001697.6
$ws.Range("K58").Value = 1131.375
$ws.Range("L58").Value = 2001697.6
$ws.Range("M58").Value = -928.375

# Row 99
$ws.Range("H99").Value = 2659.5557
$ws.Range("I99").Value = 2544.4
$ws.Range("J99").Value = 2803.5
$ws.Range("K99").Value = 2544.4
$ws.Range("L99").Value = 2803.5
$ws.Range("M99").Value = -1046.4

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Row 126
$ws.Range("H126").Value = 2659.5557
$ws.Range("I126").Value = 2544.4
$ws.Range("J126").Value = 2803.5
$ws.Range("K126").Value = 7633.200000000001
$ws.Range("L126").Value = 8410.5
$ws.Range("M126").Value = -5163.200000000001

# Row 132
$ws.Range("H132").Value = 743595.5600000001
$ws.Range("I132").Value = 479286.38
$ws.Range("J132").Value = 1668677.9
$ws.Range("K132").Value = 1437859.14
$ws.Range("L132").Value = 5006033.699999999
$ws.Range("M132").Value = -1435329.14

# Row 134
$ws.Range("H134").Value = 2489.8333
$ws.Range("I134").Value = 1882.9231
$ws.Range("J134").Value = 4067.8
$ws.Range("K134").Value = 5648.7693
$ws.Range("L134").Value = 12203.4
$ws.Range("M134").Value = -3113.7693
$ws.Range("N134").Value = -17273.4

# Row 136
$ws.Range("H136").Value = 770579.9399999999
$ws.Range("I136").Value = 1131.375
$ws.Range("J136").Value = 2001697.6
$ws.Range("K136").Value = 3394.125
$ws.Range("L136").Value = 6005092.800000001
$ws.Range("M136").Value = -844.125

$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 7999
$ws.Range("I56").Value = 7999
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 7999
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -7469

# Row 60
$ws.Range("H60").Value = 36.666668
$ws.Range("I60").Value = 36.666668
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 110.000004
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 140.999996

# Row 134
$ws.Range("H134").Value = 1804.6923
$ws.Range("I134").Value = 1804.6923
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5414.0769
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -344.0769

# Row 140
$ws.Range("H140").Value = 1734.2963
$ws.Range("I140").Value = 1245.4736
$ws.Range("J140").Value = 2895.25
$ws.Range("K140").Value = 3736.4208
$ws.Range("L140").Value = 8685.75
$ws.Range("M140").Value = 1443.5792
$ws.Range("N140").Value = -19045.75

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 411.5
$ws.Range("I12").Value = 411.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 411.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -271.5
$ws.Range("N12").ClearContents()

# Row 80
$ws.Range("H80").Value = 2598.818
$ws.Range("I80").Value = 2371.8572
$ws.Range("J80").Value = 2996
$ws.Range("K80").Value = 2371.8572
$ws.Range("L80").Value = 2996
$ws.Range("M80").Value = -1373.8572
$ws.Range("N80").Value = -4992

# Row 83
$ws.Range("H83").Value = 2598.818
$ws.Range("I83").Value = 2371.8572
$ws.Range("J83").Value = 2996
$ws.Range("K83").Value = 11859.286
$ws.Range("L83").Value = 14980
$ws.Range("M83").Value = -6867.286
$ws.Range("N83").Value = -24964

# Row 102
$ws.Range("H102").Value = 2839.6843
$ws.Range("I102").Value = 1723.091
$ws.Range("J102").Value = 4375
$ws.Range("K102").Value = 1723.091
$ws.Range("L102").Value = 4375
$ws.Range("M102").Value = -101.0909999999999
$ws.Range("N102").Value = -7619

# Row 132
$ws.Range("H132").Value = 326753.6
$ws.Range("I132").Value = 374050.8
$ws.Range("J132").Value = 7497.25
$ws.Range("K132").Value = 1122152.4
$ws.Range("L132").Value = 22491.75
$ws.Range("M132").Value = -1119622.4
$ws.Range("N132").Value = -27551.75

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 2171.7144
$ws.Range("I82").Value = 2708.5
$ws.Range("J82").Value = 1098.1428
$ws.Range("K82").Value = 2708.5
$ws.Range("L82").Value = 1098.1428
$ws.Range("M82").Value = -2347.5
$ws.Range("N82").Value = -1820.1428

# Row 85
$ws.Range("H85").Value = 2171.7144
$ws.Range("I85").Value = 2708.5
$ws.Range("J85").Value = 1098.1428
$ws.Range("K85").Value = 2708.5
$ws.Range("L85").Value = 1098.1428
$ws.Range("M85").Value = -1460.5
$ws.Range("N85").Value = -3594.1428

# Row 136
$ws.Range("H136").Value = 2569.825
$ws.Range("I136").Value = 2059.276
$ws.Range("J136").Value = 3915.818
$ws.Range("K136").Value = 6177.828
$ws.Range("L136").Value = 11747.454
$ws.Range("M136").Value = -3627.828
$ws.Range("N136").Value = -16847.454

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1857.3334
$ws.Range("I126").Value = 1545.375
$ws.Range("J126").Value = 2213.8572
$ws.Range("K126").Value = 4636.125
$ws.Range("L126").Value = 6641.571599999999
$ws.Range("M126").Value = -2166.125
$ws.Range("N126").Value = -11581.5716

# Row 132
$ws.Range("H132").Value = 300951
$ws.Range("I132").Value = 341139.78
$ws.Range("J132").Value = 4558.625
$ws.Range("K132").Value = 1023419.34
$ws.Range("L132").Value = 13675.875
$ws.Range("M132").Value = -1020889.34
$ws.Range("N132").Value = -18735.875
